$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data row is inserted at row 14 ("Region Metropolitana" entry dated
# 2021-10-18); the previously-existing rows 14-47 shift down to rows 15-48.
# Since every column in every row changes identity (not just a couple of
# cells), the full row range 14-48 is rewritten explicitly cell by cell,
# which reproduces the shift exactly and is far less error-prone than trying
# to drive Rows.Insert/EntireRow shifting semantics.

# Row 14
$ws.Cells.Item(14, 1).Value = 10
$ws.Cells.Item(14, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(14, 3).Value = "La Araucanía"
$ws.Cells.Item(14, 4).Value = 44487
$ws.Cells.Item(14, 5).Value = 9
$ws.Cells.Item(14, 6).Value = 100112022
$ws.Cells.Item(14, 7).Value = "Arveja Verde"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 110
$ws.Cells.Item(14, 11).Value = 25000
$ws.Cells.Item(14, 12).Value = 25000
$ws.Cells.Item(14, 13).Value = 25000
$ws.Cells.Item(14, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(14, 15).Value = "Región Metropolitana"
$ws.Cells.Item(14, 16).Value = 1000
$ws.Cells.Item(14, 17).Value = 25
$ws.Cells.Item(14, 18).Value = "Hortaliza"

# Row 15
$ws.Cells.Item(15, 1).Value = 10
$ws.Cells.Item(15, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(15, 3).Value = "La Araucanía"
$ws.Cells.Item(15, 4).Value = 44217
$ws.Cells.Item(15, 5).Value = 9
$ws.Cells.Item(15, 6).Value = 100112022
$ws.Cells.Item(15, 7).Value = "Arveja Verde"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 100
$ws.Cells.Item(15, 11).Value = 23000
$ws.Cells.Item(15, 12).Value = 23000
$ws.Cells.Item(15, 13).Value = 23000
$ws.Cells.Item(15, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(15, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(15, 16).Value = 920
$ws.Cells.Item(15, 17).Value = 25
$ws.Cells.Item(15, 18).Value = "Hortaliza"

# Row 16
$ws.Cells.Item(16, 1).Value = 10
$ws.Cells.Item(16, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(16, 3).Value = "La Araucanía"
$ws.Cells.Item(16, 4).Value = 44238
$ws.Cells.Item(16, 5).Value = 9
$ws.Cells.Item(16, 6).Value = 100112022
$ws.Cells.Item(16, 7).Value = "Arveja Verde"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 110
$ws.Cells.Item(16, 11).Value = 17000
$ws.Cells.Item(16, 12).Value = 17000
$ws.Cells.Item(16, 13).Value = 17000
$ws.Cells.Item(16, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(16, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(16, 16).Value = 680
$ws.Cells.Item(16, 17).Value = 25
$ws.Cells.Item(16, 18).Value = "Hortaliza"

# Row 17
$ws.Cells.Item(17, 1).Value = 10
$ws.Cells.Item(17, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(17, 3).Value = "La Araucanía"
$ws.Cells.Item(17, 4).Value = 44222
$ws.Cells.Item(17, 5).Value = 9
$ws.Cells.Item(17, 6).Value = 100112022
$ws.Cells.Item(17, 7).Value = "Arveja Verde"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 65
$ws.Cells.Item(17, 11).Value = 18000
$ws.Cells.Item(17, 12).Value = 18000
$ws.Cells.Item(17, 13).Value = 18000
$ws.Cells.Item(17, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(17, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(17, 16).Value = 720
$ws.Cells.Item(17, 17).Value = 25
$ws.Cells.Item(17, 18).Value = "Hortaliza"

# Row 18
$ws.Cells.Item(18, 1).Value = 10
$ws.Cells.Item(18, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(18, 3).Value = "La Araucanía"
$ws.Cells.Item(18, 4).Value = 44211
$ws.Cells.Item(18, 5).Value = 9
$ws.Cells.Item(18, 6).Value = 100112022
$ws.Cells.Item(18, 7).Value = "Arveja Verde"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 110
$ws.Cells.Item(18, 11).Value = 17000
$ws.Cells.Item(18, 12).Value = 17000
$ws.Cells.Item(18, 13).Value = 17000
$ws.Cells.Item(18, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(18, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(18, 16).Value = 680
$ws.Cells.Item(18, 17).Value = 25
$ws.Cells.Item(18, 18).Value = "Hortaliza"

# Row 19
$ws.Cells.Item(19, 1).Value = 10
$ws.Cells.Item(19, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(19, 3).Value = "La Araucanía"
$ws.Cells.Item(19, 4).Value = 44210
$ws.Cells.Item(19, 5).Value = 9
$ws.Cells.Item(19, 6).Value = 100112022
$ws.Cells.Item(19, 7).Value = "Arveja Verde"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 150
$ws.Cells.Item(19, 11).Value = 17000
$ws.Cells.Item(19, 12).Value = 17000
$ws.Cells.Item(19, 13).Value = 17000
$ws.Cells.Item(19, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(19, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(19, 16).Value = 680
$ws.Cells.Item(19, 17).Value = 25
$ws.Cells.Item(19, 18).Value = "Hortaliza"

# Row 20
$ws.Cells.Item(20, 1).Value = 10
$ws.Cells.Item(20, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(20, 3).Value = "La Araucanía"
$ws.Cells.Item(20, 4).Value = 44257
$ws.Cells.Item(20, 5).Value = 9
$ws.Cells.Item(20, 6).Value = 100112022
$ws.Cells.Item(20, 7).Value = "Arveja Verde"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 40
$ws.Cells.Item(20, 11).Value = 20000
$ws.Cells.Item(20, 12).Value = 20000
$ws.Cells.Item(20, 13).Value = 20000
$ws.Cells.Item(20, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(20, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(20, 16).Value = 800
$ws.Cells.Item(20, 17).Value = 25
$ws.Cells.Item(20, 18).Value = "Hortaliza"

# Row 21
$ws.Cells.Item(21, 1).Value = 10
$ws.Cells.Item(21, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(21, 3).Value = "La Araucanía"
$ws.Cells.Item(21, 4).Value = 44200
$ws.Cells.Item(21, 5).Value = 9
$ws.Cells.Item(21, 6).Value = 100112022
$ws.Cells.Item(21, 7).Value = "Arveja Verde"
$ws.Cells.Item(21, 8).Value = "Sin especificar"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 80
$ws.Cells.Item(21, 11).Value = 18000
$ws.Cells.Item(21, 12).Value = 18000
$ws.Cells.Item(21, 13).Value = 18000
$ws.Cells.Item(21, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(21, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(21, 16).Value = 720
$ws.Cells.Item(21, 17).Value = 25
$ws.Cells.Item(21, 18).Value = "Hortaliza"

# Row 22
$ws.Cells.Item(22, 1).Value = 10
$ws.Cells.Item(22, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(22, 3).Value = "La Araucanía"
$ws.Cells.Item(22, 4).Value = 44225
$ws.Cells.Item(22, 5).Value = 9
$ws.Cells.Item(22, 6).Value = 100112022
$ws.Cells.Item(22, 7).Value = "Arveja Verde"
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 25
$ws.Cells.Item(22, 11).Value = 17000
$ws.Cells.Item(22, 12).Value = 17000
$ws.Cells.Item(22, 13).Value = 17000
$ws.Cells.Item(22, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(22, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(22, 16).Value = 680
$ws.Cells.Item(22, 17).Value = 25
$ws.Cells.Item(22, 18).Value = "Hortaliza"

# Row 23
$ws.Cells.Item(23, 1).Value = 10
$ws.Cells.Item(23, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(23, 3).Value = "La Araucanía"
$ws.Cells.Item(23, 4).Value = 44175
$ws.Cells.Item(23, 5).Value = 9
$ws.Cells.Item(23, 6).Value = 100112022
$ws.Cells.Item(23, 7).Value = "Arveja Verde"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 40
$ws.Cells.Item(23, 11).Value = 20000
$ws.Cells.Item(23, 12).Value = 20000
$ws.Cells.Item(23, 13).Value = 20000
$ws.Cells.Item(23, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(23, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(23, 16).Value = 800
$ws.Cells.Item(23, 17).Value = 25
$ws.Cells.Item(23, 18).Value = "Hortaliza"

# Row 24
$ws.Cells.Item(24, 1).Value = 10
$ws.Cells.Item(24, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(24, 3).Value = "La Araucanía"
$ws.Cells.Item(24, 4).Value = 44188
$ws.Cells.Item(24, 5).Value = 9
$ws.Cells.Item(24, 6).Value = 100112022
$ws.Cells.Item(24, 7).Value = "Arveja Verde"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 50
$ws.Cells.Item(24, 11).Value = 25000
$ws.Cells.Item(24, 12).Value = 25000
$ws.Cells.Item(24, 13).Value = 25000
$ws.Cells.Item(24, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(24, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(24, 16).Value = 1000
$ws.Cells.Item(24, 17).Value = 25
$ws.Cells.Item(24, 18).Value = "Hortaliza"

# Row 25
$ws.Cells.Item(25, 1).Value = 10
$ws.Cells.Item(25, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(25, 3).Value = "La Araucanía"
$ws.Cells.Item(25, 4).Value = 44271
$ws.Cells.Item(25, 5).Value = 9
$ws.Cells.Item(25, 6).Value = 100112022
$ws.Cells.Item(25, 7).Value = "Arveja Verde"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 80
$ws.Cells.Item(25, 11).Value = 18000
$ws.Cells.Item(25, 12).Value = 18000
$ws.Cells.Item(25, 13).Value = 18000
$ws.Cells.Item(25, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(25, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(25, 16).Value = 720
$ws.Cells.Item(25, 17).Value = 25
$ws.Cells.Item(25, 18).Value = "Hortaliza"

# Row 26
$ws.Cells.Item(26, 1).Value = 10
$ws.Cells.Item(26, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(26, 3).Value = "La Araucanía"
$ws.Cells.Item(26, 4).Value = 44483
$ws.Cells.Item(26, 5).Value = 9
$ws.Cells.Item(26, 6).Value = 100112022
$ws.Cells.Item(26, 7).Value = "Arveja Verde"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 20
$ws.Cells.Item(26, 11).Value = 25000
$ws.Cells.Item(26, 12).Value = 25000
$ws.Cells.Item(26, 13).Value = 25000
$ws.Cells.Item(26, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(26, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(26, 16).Value = 1000
$ws.Cells.Item(26, 17).Value = 25
$ws.Cells.Item(26, 18).Value = "Hortaliza"

# Row 27
$ws.Cells.Item(27, 1).Value = 10
$ws.Cells.Item(27, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(27, 3).Value = "La Araucanía"
$ws.Cells.Item(27, 4).Value = 44224
$ws.Cells.Item(27, 5).Value = 9
$ws.Cells.Item(27, 6).Value = 100112022
$ws.Cells.Item(27, 7).Value = "Arveja Verde"
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 320
$ws.Cells.Item(27, 11).Value = 17000
$ws.Cells.Item(27, 12).Value = 17000
$ws.Cells.Item(27, 13).Value = 17000
$ws.Cells.Item(27, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(27, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(27, 16).Value = 680
$ws.Cells.Item(27, 17).Value = 25
$ws.Cells.Item(27, 18).Value = "Hortaliza"

# Row 28
$ws.Cells.Item(28, 1).Value = 10
$ws.Cells.Item(28, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(28, 3).Value = "La Araucanía"
$ws.Cells.Item(28, 4).Value = 44195
$ws.Cells.Item(28, 5).Value = 9
$ws.Cells.Item(28, 6).Value = 100112022
$ws.Cells.Item(28, 7).Value = "Arveja Verde"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 110
$ws.Cells.Item(28, 11).Value = 25000
$ws.Cells.Item(28, 12).Value = 25000
$ws.Cells.Item(28, 13).Value = 25000
$ws.Cells.Item(28, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(28, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(28, 16).Value = 1000
$ws.Cells.Item(28, 17).Value = 25
$ws.Cells.Item(28, 18).Value = "Hortaliza"

# Row 29
$ws.Cells.Item(29, 1).Value = 10
$ws.Cells.Item(29, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(29, 3).Value = "La Araucanía"
$ws.Cells.Item(29, 4).Value = 44484
$ws.Cells.Item(29, 5).Value = 9
$ws.Cells.Item(29, 6).Value = 100112022
$ws.Cells.Item(29, 7).Value = "Arveja Verde"
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 10
$ws.Cells.Item(29, 11).Value = 25000
$ws.Cells.Item(29, 12).Value = 25000
$ws.Cells.Item(29, 13).Value = 25000
$ws.Cells.Item(29, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(29, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(29, 16).Value = 1000
$ws.Cells.Item(29, 17).Value = 25
$ws.Cells.Item(29, 18).Value = "Hortaliza"

# Row 30
$ws.Cells.Item(30, 1).Value = 10
$ws.Cells.Item(30, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(30, 3).Value = "La Araucanía"
$ws.Cells.Item(30, 4).Value = 44232
$ws.Cells.Item(30, 5).Value = 9
$ws.Cells.Item(30, 6).Value = 100112022
$ws.Cells.Item(30, 7).Value = "Arveja Verde"
$ws.Cells.Item(30, 8).Value = "Sin especificar"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 110
$ws.Cells.Item(30, 11).Value = 20000
$ws.Cells.Item(30, 12).Value = 20000
$ws.Cells.Item(30, 13).Value = 20000
$ws.Cells.Item(30, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(30, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(30, 16).Value = 800
$ws.Cells.Item(30, 17).Value = 25
$ws.Cells.Item(30, 18).Value = "Hortaliza"

# Row 31
$ws.Cells.Item(31, 1).Value = 10
$ws.Cells.Item(31, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(31, 3).Value = "La Araucanía"
$ws.Cells.Item(31, 4).Value = 44169
$ws.Cells.Item(31, 5).Value = 9
$ws.Cells.Item(31, 6).Value = 100112022
$ws.Cells.Item(31, 7).Value = "Arveja Verde"
$ws.Cells.Item(31, 8).Value = "Sin especificar"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 135
$ws.Cells.Item(31, 11).Value = 20000
$ws.Cells.Item(31, 12).Value = 22000
$ws.Cells.Item(31, 13).Value = 20593
$ws.Cells.Item(31, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(31, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(31, 16).Value = 824
$ws.Cells.Item(31, 17).Value = 25
$ws.Cells.Item(31, 18).Value = "Hortaliza"

# Row 32
$ws.Cells.Item(32, 1).Value = 10
$ws.Cells.Item(32, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(32, 3).Value = "La Araucanía"
$ws.Cells.Item(32, 4).Value = 44273
$ws.Cells.Item(32, 5).Value = 9
$ws.Cells.Item(32, 6).Value = 100112022
$ws.Cells.Item(32, 7).Value = "Arveja Verde"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 50
$ws.Cells.Item(32, 11).Value = 20000
$ws.Cells.Item(32, 12).Value = 20000
$ws.Cells.Item(32, 13).Value = 20000
$ws.Cells.Item(32, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(32, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(32, 16).Value = 800
$ws.Cells.Item(32, 17).Value = 25
$ws.Cells.Item(32, 18).Value = "Hortaliza"

# Row 33
$ws.Cells.Item(33, 1).Value = 10
$ws.Cells.Item(33, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(33, 3).Value = "La Araucanía"
$ws.Cells.Item(33, 4).Value = 44236
$ws.Cells.Item(33, 5).Value = 9
$ws.Cells.Item(33, 6).Value = 100112022
$ws.Cells.Item(33, 7).Value = "Arveja Verde"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 95
$ws.Cells.Item(33, 11).Value = 17000
$ws.Cells.Item(33, 12).Value = 17000
$ws.Cells.Item(33, 13).Value = 17000
$ws.Cells.Item(33, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(33, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(33, 16).Value = 680
$ws.Cells.Item(33, 17).Value = 25
$ws.Cells.Item(33, 18).Value = "Hortaliza"

# Row 34
$ws.Cells.Item(34, 1).Value = 10
$ws.Cells.Item(34, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(34, 3).Value = "La Araucanía"
$ws.Cells.Item(34, 4).Value = 44235
$ws.Cells.Item(34, 5).Value = 9
$ws.Cells.Item(34, 6).Value = 100112022
$ws.Cells.Item(34, 7).Value = "Arveja Verde"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 250
$ws.Cells.Item(34, 11).Value = 17000
$ws.Cells.Item(34, 12).Value = 17000
$ws.Cells.Item(34, 13).Value = 17000
$ws.Cells.Item(34, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(34, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(34, 16).Value = 680
$ws.Cells.Item(34, 17).Value = 25
$ws.Cells.Item(34, 18).Value = "Hortaliza"

# Row 35
$ws.Cells.Item(35, 1).Value = 10
$ws.Cells.Item(35, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(35, 3).Value = "La Araucanía"
$ws.Cells.Item(35, 4).Value = 44203
$ws.Cells.Item(35, 5).Value = 9
$ws.Cells.Item(35, 6).Value = 100112022
$ws.Cells.Item(35, 7).Value = "Arveja Verde"
$ws.Cells.Item(35, 8).Value = "Sin especificar"
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 160
$ws.Cells.Item(35, 11).Value = 15000
$ws.Cells.Item(35, 12).Value = 18000
$ws.Cells.Item(35, 13).Value = 16875
$ws.Cells.Item(35, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(35, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(35, 16).Value = 675
$ws.Cells.Item(35, 17).Value = 25
$ws.Cells.Item(35, 18).Value = "Hortaliza"

# Row 36
$ws.Cells.Item(36, 1).Value = 10
$ws.Cells.Item(36, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(36, 3).Value = "La Araucanía"
$ws.Cells.Item(36, 4).Value = 44168
$ws.Cells.Item(36, 5).Value = 9
$ws.Cells.Item(36, 6).Value = 100112022
$ws.Cells.Item(36, 7).Value = "Arveja Verde"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 75
$ws.Cells.Item(36, 11).Value = 20000
$ws.Cells.Item(36, 12).Value = 20000
$ws.Cells.Item(36, 13).Value = 20000
$ws.Cells.Item(36, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(36, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(36, 16).Value = 800
$ws.Cells.Item(36, 17).Value = 25
$ws.Cells.Item(36, 18).Value = "Hortaliza"

# Row 37
$ws.Cells.Item(37, 1).Value = 10
$ws.Cells.Item(37, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(37, 3).Value = "La Araucanía"
$ws.Cells.Item(37, 4).Value = 44161
$ws.Cells.Item(37, 5).Value = 9
$ws.Cells.Item(37, 6).Value = 100112022
$ws.Cells.Item(37, 7).Value = "Arveja Verde"
$ws.Cells.Item(37, 8).Value = "Sin especificar"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 100
$ws.Cells.Item(37, 11).Value = 23000
$ws.Cells.Item(37, 12).Value = 23000
$ws.Cells.Item(37, 13).Value = 23000
$ws.Cells.Item(37, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(37, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(37, 16).Value = 920
$ws.Cells.Item(37, 17).Value = 25
$ws.Cells.Item(37, 18).Value = "Hortaliza"

# Row 38
$ws.Cells.Item(38, 1).Value = 10
$ws.Cells.Item(38, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(38, 3).Value = "La Araucanía"
$ws.Cells.Item(38, 4).Value = 44161
$ws.Cells.Item(38, 5).Value = 9
$ws.Cells.Item(38, 6).Value = 100112022
$ws.Cells.Item(38, 7).Value = "Arveja Verde"
$ws.Cells.Item(38, 8).Value = "Sin especificar"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 100
$ws.Cells.Item(38, 11).Value = 20000
$ws.Cells.Item(38, 12).Value = 21000
$ws.Cells.Item(38, 13).Value = 20500
$ws.Cells.Item(38, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(38, 15).Value = "Región del Maule"
$ws.Cells.Item(38, 16).Value = 820
$ws.Cells.Item(38, 17).Value = 25
$ws.Cells.Item(38, 18).Value = "Hortaliza"

# Row 39
$ws.Cells.Item(39, 1).Value = 10
$ws.Cells.Item(39, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(39, 3).Value = "La Araucanía"
$ws.Cells.Item(39, 4).Value = 44165
$ws.Cells.Item(39, 5).Value = 9
$ws.Cells.Item(39, 6).Value = 100112022
$ws.Cells.Item(39, 7).Value = "Arveja Verde"
$ws.Cells.Item(39, 8).Value = "Sin especificar"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 200
$ws.Cells.Item(39, 11).Value = 22000
$ws.Cells.Item(39, 12).Value = 22000
$ws.Cells.Item(39, 13).Value = 22000
$ws.Cells.Item(39, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(39, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(39, 16).Value = 880
$ws.Cells.Item(39, 17).Value = 25
$ws.Cells.Item(39, 18).Value = "Hortaliza"

# Row 40
$ws.Cells.Item(40, 1).Value = 10
$ws.Cells.Item(40, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(40, 3).Value = "La Araucanía"
$ws.Cells.Item(40, 4).Value = 44231
$ws.Cells.Item(40, 5).Value = 9
$ws.Cells.Item(40, 6).Value = 100112022
$ws.Cells.Item(40, 7).Value = "Arveja Verde"
$ws.Cells.Item(40, 8).Value = "Perfection"
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 110
$ws.Cells.Item(40, 11).Value = 20000
$ws.Cells.Item(40, 12).Value = 20000
$ws.Cells.Item(40, 13).Value = 20000
$ws.Cells.Item(40, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(40, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(40, 16).Value = 800
$ws.Cells.Item(40, 17).Value = 25
$ws.Cells.Item(40, 18).Value = "Hortaliza"

# Row 41
$ws.Cells.Item(41, 1).Value = 10
$ws.Cells.Item(41, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(41, 3).Value = "La Araucanía"
$ws.Cells.Item(41, 4).Value = 44204
$ws.Cells.Item(41, 5).Value = 9
$ws.Cells.Item(41, 6).Value = 100112022
$ws.Cells.Item(41, 7).Value = "Arveja Verde"
$ws.Cells.Item(41, 8).Value = "Sin especificar"
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 40
$ws.Cells.Item(41, 11).Value = 15000
$ws.Cells.Item(41, 12).Value = 15000
$ws.Cells.Item(41, 13).Value = 15000
$ws.Cells.Item(41, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(41, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(41, 16).Value = 600
$ws.Cells.Item(41, 17).Value = 25
$ws.Cells.Item(41, 18).Value = "Hortaliza"

# Row 42
$ws.Cells.Item(42, 1).Value = 10
$ws.Cells.Item(42, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(42, 3).Value = "La Araucanía"
$ws.Cells.Item(42, 4).Value = 44186
$ws.Cells.Item(42, 5).Value = 9
$ws.Cells.Item(42, 6).Value = 100112022
$ws.Cells.Item(42, 7).Value = "Arveja Verde"
$ws.Cells.Item(42, 8).Value = "Sin especificar"
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 40
$ws.Cells.Item(42, 11).Value = 20000
$ws.Cells.Item(42, 12).Value = 20000
$ws.Cells.Item(42, 13).Value = 20000
$ws.Cells.Item(42, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(42, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(42, 16).Value = 800
$ws.Cells.Item(42, 17).Value = 25
$ws.Cells.Item(42, 18).Value = "Hortaliza"

# Row 43
$ws.Cells.Item(43, 1).Value = 10
$ws.Cells.Item(43, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(43, 3).Value = "La Araucanía"
$ws.Cells.Item(43, 4).Value = 44243
$ws.Cells.Item(43, 5).Value = 9
$ws.Cells.Item(43, 6).Value = 100112022
$ws.Cells.Item(43, 7).Value = "Arveja Verde"
$ws.Cells.Item(43, 8).Value = "Sin especificar"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 95
$ws.Cells.Item(43, 11).Value = 17000
$ws.Cells.Item(43, 12).Value = 17000
$ws.Cells.Item(43, 13).Value = 17000
$ws.Cells.Item(43, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(43, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(43, 16).Value = 680
$ws.Cells.Item(43, 17).Value = 25
$ws.Cells.Item(43, 18).Value = "Hortaliza"

# Row 44
$ws.Cells.Item(44, 1).Value = 10
$ws.Cells.Item(44, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(44, 3).Value = "La Araucanía"
$ws.Cells.Item(44, 4).Value = 44242
$ws.Cells.Item(44, 5).Value = 9
$ws.Cells.Item(44, 6).Value = 100112022
$ws.Cells.Item(44, 7).Value = "Arveja Verde"
$ws.Cells.Item(44, 8).Value = "Sin especificar"
$ws.Cells.Item(44, 9).Value = "Primera"
$ws.Cells.Item(44, 10).Value = 85
$ws.Cells.Item(44, 11).Value = 17000
$ws.Cells.Item(44, 12).Value = 17000
$ws.Cells.Item(44, 13).Value = 17000
$ws.Cells.Item(44, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(44, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(44, 16).Value = 680
$ws.Cells.Item(44, 17).Value = 25
$ws.Cells.Item(44, 18).Value = "Hortaliza"

# Row 45
$ws.Cells.Item(45, 1).Value = 10
$ws.Cells.Item(45, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(45, 3).Value = "La Araucanía"
$ws.Cells.Item(45, 4).Value = 44159
$ws.Cells.Item(45, 5).Value = 9
$ws.Cells.Item(45, 6).Value = 100112022
$ws.Cells.Item(45, 7).Value = "Arveja Verde"
$ws.Cells.Item(45, 8).Value = "Sin especificar"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 100
$ws.Cells.Item(45, 11).Value = 20000
$ws.Cells.Item(45, 12).Value = 20000
$ws.Cells.Item(45, 13).Value = 20000
$ws.Cells.Item(45, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(45, 15).Value = "Región del Maule"
$ws.Cells.Item(45, 16).Value = 800
$ws.Cells.Item(45, 17).Value = 25
$ws.Cells.Item(45, 18).Value = "Hortaliza"

# Row 46
$ws.Cells.Item(46, 1).Value = 10
$ws.Cells.Item(46, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(46, 3).Value = "La Araucanía"
$ws.Cells.Item(46, 4).Value = 44166
$ws.Cells.Item(46, 5).Value = 9
$ws.Cells.Item(46, 6).Value = 100112022
$ws.Cells.Item(46, 7).Value = "Arveja Verde"
$ws.Cells.Item(46, 8).Value = "Sin especificar"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 100
$ws.Cells.Item(46, 11).Value = 21000
$ws.Cells.Item(46, 12).Value = 22000
$ws.Cells.Item(46, 13).Value = 21450
$ws.Cells.Item(46, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(46, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(46, 16).Value = 858
$ws.Cells.Item(46, 17).Value = 25
$ws.Cells.Item(46, 18).Value = "Hortaliza"

# Row 47
$ws.Cells.Item(47, 1).Value = 10
$ws.Cells.Item(47, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(47, 3).Value = "La Araucanía"
$ws.Cells.Item(47, 4).Value = 44172
$ws.Cells.Item(47, 5).Value = 9
$ws.Cells.Item(47, 6).Value = 100112022
$ws.Cells.Item(47, 7).Value = "Arveja Verde"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 85
$ws.Cells.Item(47, 11).Value = 19000
$ws.Cells.Item(47, 12).Value = 20000
$ws.Cells.Item(47, 13).Value = 19412
$ws.Cells.Item(47, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(47, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(47, 16).Value = 776
$ws.Cells.Item(47, 17).Value = 25
$ws.Cells.Item(47, 18).Value = "Hortaliza"

# Row 48
$ws.Cells.Item(48, 1).Value = 10
$ws.Cells.Item(48, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(48, 3).Value = "La Araucanía"
$ws.Cells.Item(48, 4).Value = 44201
$ws.Cells.Item(48, 5).Value = 9
$ws.Cells.Item(48, 6).Value = 100112022
$ws.Cells.Item(48, 7).Value = "Arveja Verde"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 30
$ws.Cells.Item(48, 11).Value = 18000
$ws.Cells.Item(48, 12).Value = 18000
$ws.Cells.Item(48, 13).Value = 18000
$ws.Cells.Item(48, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(48, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(48, 16).Value = 720
$ws.Cells.Item(48, 17).Value = 25
$ws.Cells.Item(48, 18).Value = "Hortaliza"

# Row 48 is a newly-created row (the sheet previously ended at row 47), so its
# date cell needs the same date number format the rest of column D already has.
$ws.Cells.Item(48, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
